# Apply the changes described by the diff:
#  - add a new row 11 with a new "УИК №142" record that has Null-values in LineNumber
#  - add two new shared strings ("УИК №142", "Null-values") implicitly via the new cell values
#  - change active selection from A1048576 to A11 (and scroll the top-left cell back to A1)
#  - set explicit per-column widths instead of a single uniform width for all columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 11), mirroring the other UIK rows in the sheet ---
$ws.Range("A11").Value = "Республика Дагестан"
$ws.Range("B11").Value = "10 Республика Дагестан - Северный"
$ws.Range("C11").Value = "5 Бабаюртовская"
$ws.Range("D11").Value = "УИК №142"
$ws.Range("AK11").Value = "Null-values"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth  = 19.1666666666667
$ws.Columns.Item(2).ColumnWidth  = 32
$ws.Columns.Item(3).ColumnWidth  = 15.5
$ws.Columns.Item(4).ColumnWidth  = 10
$ws.Range($ws.Columns.Item(5), $ws.Columns.Item(6)).ColumnWidth = 5
$ws.Columns.Item(7).ColumnWidth  = 1.83333333333333
$ws.Columns.Item(8).ColumnWidth  = 5
$ws.Columns.Item(9).ColumnWidth  = 2.83333333333333
$ws.Columns.Item(10).ColumnWidth = 4
$ws.Columns.Item(11).ColumnWidth = 2.83333333333333
$ws.Columns.Item(12).ColumnWidth = 5
$ws.Columns.Item(13).ColumnWidth = 2.83333333333333
$ws.Columns.Item(14).ColumnWidth = 5
$ws.Range($ws.Columns.Item(15), $ws.Columns.Item(25)).ColumnWidth = 2.83333333333333
$ws.Columns.Item(26).ColumnWidth = 5
$ws.Range($ws.Columns.Item(27), $ws.Columns.Item(36)).ColumnWidth = 2.83333333333333
$ws.Columns.Item(37).ColumnWidth = 192.333333333333

# --- View: move back to top of sheet and select the newly added row ---
$null = $ws.Range("A1").Select()
$null = $ws.Range("A11").Select()
